$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. "synapse organization" GO term gains a forced. prefix
$ws.Range("C6").Value = "forced.GO:0050808"

# 2. The gene list that used to sit on row 7 (D7:Y7, alongside
#    "modulation of chemical synaptic transmission" / GO:0050804) moves
#    down to row 10, landing next to "regulation of angiogenesis", whose
#    GO id also gains a forced. prefix. Rows 8/9 (chemical synaptic
#    transmission / positive regulation of vasculature development) shift
#    up to take rows 7/8... actually they keep their row numbers; the
#    gene list vacates row 7 and reappears on row 10.
$genes = $ws.Range("D7:Y7").Value()

$ws.Range("D10:Y10").Value = $genes
$ws.Range("C10").Value = "forced.GO:0045765"

# Row 7 only keeps its first three cells now (astrocytes / modulation of
# chemical synaptic transmission / GO:0050804); clear the gene list cells
# that were copied away.
$ws.Range("D7:Y7").ClearContents()

# Restore the selection state recorded in the workbook.
$ws.Range("A4").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("C8").Select()
